$d = $word.ActiveDocument

# --- 1. Turn the first blank paragraph after the AHU-mode statistics into a
#        new "List Bullet" item: "No faults were found ..." ----------------
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*mechanical cooling mode: 26.48%*") {
        $anchorPara = $p
    }
}
$anchorIndex = $anchorPara.Index

# Inserting right after an existing "List Bullet" paragraph makes the brand
# new paragraph inherit that same style automatically, which reproduces the
# target <w:pPr><w:pStyle w:val="ListBullet"/></w:pPr> markup without
# stamping an extraneous rsid attribute onto the paragraph mark.
$anchorPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($anchorIndex + 1)
$newPara.Range.Text = "No faults were found in this given dataset for the equation defined by ASHRAE."

# InsertParagraphAfter() turned the original two blank paragraphs into three;
# delete the extra one so a single blank paragraph remains in between, just
# like in the source document.
$d.Paragraphs($anchorIndex + 3).Range.Delete()

# --- 2. Reword the control-system-tuning suggestion -------------------------
$d.Content.Find.Execute(
    "The max value found is 0 changes per hour which is low or does not appear to be an issue requiring control system tuning.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "No control system tuning appears to be needed for the operating conditions of this AHU.",
    2)

# --- 3. Bump the "Report generated" timestamp -------------------------------
$d.Content.Find.Execute(
    "Report generated: Wed Feb 15 11:22:09 2023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Report generated: Thu Feb 16 09:16:31 2023",
    2)
